$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (old F "District" shifts to G)
$ws.Columns.Item(6).Insert()

# New header for inserted column F
$ws.Range("F2").Value = "Address"

# Populate the new Address column (F) for each data row.
$ws.Range("F3").Value = "Sri Mahantesh High School Kamalapur"
$ws.Range("F4").Value = "G H S BelagaliHubli"
$ws.Range("F5").Value = "S P H S Honnapur"
$ws.Range("F6").Value = "G H S UgginakeriKalaghatagi"
$ws.Range("F7").Value = "Sree Y V M H School MorabNavalgund"
$ws.Range("F8").Value = "Govt. High School Shivapur Kalaghatagi"
$ws.Range("F9").Value = "G H S Gandhinagar"
$ws.Range("F10").Value = "G H S ShisvinahalliNavalgund"
$ws.Range("F11").Value = "G H S NalawadiNavalgund"
$ws.Range("F12").Value = "G H S AdaragunchiHubballi"
$ws.Range("F13").Value = "Nehru Comp. PU CollegeHebballi"
$ws.Range("F14").Value = "G H S Shivalli"
$ws.Range("F15").Value = "G H S Devarahubballi"
$ws.Range("F16").Value = "G H S Maradagi"
$ws.Range("F17").Value = "Govt. High SchoolDevikoppaKalaghatagi"
$ws.Range("F18").Value = "Govt. High SchoolByalyalNavalgund"
$ws.Range("F19").Value = "G H S Narendra"
# F20 intentionally left blank - no address text is available for this row
# (the source name/address cell for row 20 only contains a single segment,
# which is treated as the district rather than a separate street address).
$ws.Range("F21").Value = "Pandit Nehru High School ShiraguppiHubli"
$ws.Range("F22").Value = "H M H S Mugad"
$ws.Range("F23").Value = "G H S ShirurNavalgund"
$ws.Range("F24").Value = "S J A N High School Navalgund"
$ws.Range("F25").Value = "Shri Durgadevi High School Dajibanpeth Hubballi"
$ws.Range("F26").Value = "G H S KuruvinakoppaKalaghatagi"
